$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'45.967.32"
$ws.Range("E2").Value = "'  -0.22%  "

$ws.Range("D3").Value = "'2.611.24"

$ws.Range("E4").Value = "'  -0.05%  "

$ws.Range("D5").Value = "'308.91"
$ws.Range("E5").Value = "'  +1.07%  "

$ws.Range("D6").Value = "'99.04"
$ws.Range("E6").Value = "'  -0.67%  "

$ws.Range("D7").Value = "'0.595"

$ws.Range("E8").Value = "'  -0.02%  "

$ws.Range("D9").Value = "'0.580"
$ws.Range("E9").Value = "'  +1.28%  "

$ws.Range("D10").Value = "'38.84"
$ws.Range("E10").Value = "'  +1.03%  "

$ws.Range("D11").Value = "'0.0843"
$ws.Range("E11").Value = "'  +0.79%  "

$ws.Range("D12").Value = "'54.20"
$ws.Range("E12").Value = "'  -0.68%  "

$ws.Range("D13").Value = "'8.09"
$ws.Range("E13").Value = "'  -1.77%  "

$ws.Range("D14").Value = "'3.018.73"
$ws.Range("E14").Value = "'  +1.26%  "

$ws.Range("D16").Value = "'2.614.70"
$ws.Range("E16").Value = "'  +0.72%  "

$ws.Range("D17").Value = "'0.913"
$ws.Range("E17").Value = "'  +1.47%  "

$ws.Range("D18").Value = "'14.82"
$ws.Range("E18").Value = "'  +0.25%  "

$ws.Range("D19").Value = "'46.269.40"
$ws.Range("E19").Value = "'  +0.03%  "

$ws.Range("E20").Value = "'  +0.84%  "

$ws.Range("D21").Value = "'6.74"
$ws.Range("E21").Value = "'  +1.81%  "

$ws.Range("D22").Value = "'12.71"
$ws.Range("E22").Value = "'  -1.74%  "

$ws.Range("D23").Value = "'291.70"
$ws.Range("E23").Value = "'  +15.44%  "

$ws.Range("D24").Value = "'72.70"
$ws.Range("E24").Value = "'  +2.31%  "

$ws.Range("D25").Value = "'3.04"
$ws.Range("E25").Value = "'  +1.91%  "

$ws.Range("D26").Value = "'2.24"
$ws.Range("E26").Value = "'  +2.93%  "

$ws.Range("D27").Value = "'29.68"
$ws.Range("E27").Value = "'  +5.85%  "

$ws.Range("E28").Value = "'  -0.03%  "

$ws.Range("D29").Value = "'4.05"
$ws.Range("E29").Value = "'  +0.98%  "

$ws.Range("D30").Value = "'10.76"
$ws.Range("E30").Value = "'  +3.43%  "

$ws.Range("D31").Value = "'38.76"
$ws.Range("E31").Value = "'  -1.97%  "

$ws.Range("D32").Value = "'2.21"
$ws.Range("E32").Value = "'  -2.52%  "

$ws.Range("D33").Value = "'6.26"
$ws.Range("E33").Value = "'  +3.69%  "

$ws.Range("D34").Value = "'160.43"
$ws.Range("E34").Value = "'  +5.20%  "

$ws.Range("D35").Value = "'3.62"
$ws.Range("E35").Value = "'  -1.21%  "

$ws.Range("E36").Value = "'  -1.56%  "

$ws.Range("D37").Value = "'0.0842"
$ws.Range("E37").Value = "'  +2.35%  "

$ws.Range("E38").Value = "'  -3.93%  "

$ws.Range("E39").Value = "'  +4.67%  "

$ws.Range("D40").Value = "'0.123"
$ws.Range("E40").Value = "'  +1.39%  "

$ws.Range("D41").Value = "'15.69"
$ws.Range("E41").Value = "'  -2.53%  "

$ws.Range("E42").Value = "'  +3.21%  "

$ws.Range("D43").Value = "'3.55"
$ws.Range("E43").Value = "'  -1.26%  "

$ws.Range("D44").Value = "'21.57"
$ws.Range("E44").Value = "'  +8.13%  "

$ws.Range("D45").Value = "'4.01"
$ws.Range("E45").Value = "'  -3.72%  "

$ws.Range("D46").Value = "'2.114.61"
$ws.Range("E46").Value = "'  +2.87%  "

$ws.Range("D47").Value = "'96.09"
$ws.Range("E47").Value = "'  +5.79%  "

$ws.Range("E48").Value = "'  -0.17%  "

$ws.Range("D49").Value = "'9.38"
$ws.Range("E49").Value = "'  +1.26%  "

$ws.Range("D50").Value = "'109.60"
$ws.Range("E50").Value = "'  +1.33%  "

$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "'2.870.49"
$ws.Range("E51").Value = "'  +1.01%  "
